# Slide 2, "Content Placeholder 2" shape, 3rd paragraph:
#   "Allows users to visualize and share data (sensors values, media, web links etc.)"
# is split into three runs:
#   "Allows users to visualize and "
#   "anonymously share data and digital content "
#   "(sensors values, media, web links etc.)"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$para = $tr.Paragraphs(3, 1)

$prefix = "Allows users to visualize and "
$oldMiddle = "share data "
$newMiddle = "anonymously share data and digital content "

# Locate the "share data " substring (right after the unchanged prefix) inside
# the paragraph and replace it in place so the run splits line up with the
# original wording that stays untouched before/after it.
$middleStart = $para.Start + $prefix.Length
$middleRange = $tr.Characters($middleStart, $oldMiddle.Length)
$middleRange.Text = $newMiddle
